# Add new date column (07_02_2024) with receptionist appointment-status-change
# counts for the week of Jan 29 - Feb 7, 2024.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column G
$ws.Range("G1").Value = "07_02_2024"

# Values for each receptionist row
$ws.Range("G2").Value = 3063
$ws.Range("G3").Value = 2314
$ws.Range("G4").Value = 3618
$ws.Range("G5").Value = 6927

# Update the active cell/selection to reflect where the user ended up
$ws.Range("I9").Select()
